# vg fund risk metrics
# Update a handful of slightly-recomputed "Adj Close" (F) values and
# append five new trading-day rows (2023-10-31 .. 2023-11-06) to the
# Yahoo Finance historical-data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-computed Adj Close values for existing rows ---------------------
$ws.Range("F7").Value = 84.340622
$ws.Range("F8").Value = 86.238838
$ws.Range("F10").Value = 87.636009
$ws.Range("F12").Value = 89.44750999999999
$ws.Range("F13").Value = 88.77301
$ws.Range("F15").Value = 86.238838
$ws.Range("F17").Value = 88.22378500000001
$ws.Range("F19").Value = 87.664917
$ws.Range("F21").Value = 87.934715
$ws.Range("F22").Value = 89.534233
$ws.Range("F24").Value = 88.08889000000001
$ws.Range("F28").Value = 85.602898
$ws.Range("F34").Value = 83.155449
$ws.Range("F36").Value = 82.19188699999999
$ws.Range("F44").Value = 82.403862
$ws.Range("F45").Value = 83.11689800000001
$ws.Range("F47").Value = 83.290352
$ws.Range("F50").Value = 85.082573
$ws.Range("F51").Value = 87.154228
$ws.Range("F53").Value = 87.99252300000001
$ws.Range("F54").Value = 86.82662999999999
$ws.Range("F60").Value = 85.082573
$ws.Range("F62").Value = 85.381271

# --- Append new historical rows 210-214 ----------------------------------
$ws.Cells.Item(210, 1).NumberFormat = "@"
$ws.Cells.Item(210, 1).Value = "2023-10-31"
$ws.Cells.Item(210, 2).Value = 64.519997
$ws.Cells.Item(210, 3).Value = 64.970001
$ws.Cells.Item(210, 4).Value = 63.549999
$ws.Cells.Item(210, 5).Value = 63.740002
$ws.Cells.Item(210, 6).Value = 63.740002
$ws.Cells.Item(210, 7).Value = 501700

$ws.Cells.Item(211, 1).NumberFormat = "@"
$ws.Cells.Item(211, 1).Value = "2023-11-01"
$ws.Cells.Item(211, 2).Value = 64.400002
$ws.Cells.Item(211, 3).Value = 65.620003
$ws.Cells.Item(211, 4).Value = 64.400002
$ws.Cells.Item(211, 5).Value = 65.41999800000001
$ws.Cells.Item(211, 6).Value = 65.41999800000001
$ws.Cells.Item(211, 7).Value = 611400

$ws.Cells.Item(212, 1).NumberFormat = "@"
$ws.Cells.Item(212, 1).Value = "2023-11-02"
$ws.Cells.Item(212, 2).Value = 67.33000199999999
$ws.Cells.Item(212, 3).Value = 68.120003
$ws.Cells.Item(212, 4).Value = 66.870003
$ws.Cells.Item(212, 5).Value = 67.800003
$ws.Cells.Item(212, 6).Value = 67.800003
$ws.Cells.Item(212, 7).Value = 737200

$ws.Cells.Item(213, 1).NumberFormat = "@"
$ws.Cells.Item(213, 1).Value = "2023-11-03"
$ws.Cells.Item(213, 2).Value = 69.81999999999999
$ws.Cells.Item(213, 3).Value = 69.959999
$ws.Cells.Item(213, 4).Value = 68.300003
$ws.Cells.Item(213, 5).Value = 68.370003
$ws.Cells.Item(213, 6).Value = 68.370003
$ws.Cells.Item(213, 7).Value = 821900

$ws.Cells.Item(214, 1).NumberFormat = "@"
$ws.Cells.Item(214, 1).Value = "2023-11-06"
$ws.Cells.Item(214, 2).Value = 67.800003
$ws.Cells.Item(214, 3).Value = 67.94000200000001
$ws.Cells.Item(214, 4).Value = 67.290001
$ws.Cells.Item(214, 5).Value = 67.519997
$ws.Cells.Item(214, 6).Value = 67.519997
$ws.Cells.Item(214, 7).Value = 500742

